$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new text value, derived from the commit diff (coinranking
# price/volume refresh + a block of rows that shifted down by one position).
$updates = @(
    @('D2', '244.14'),
    @('E2', '-0.40%'),
    @('D3', '26.41'),
    @('E3', '3.47%'),
    @('D4', '5.157'),
    @('E4', '0.85%'),
    @('D5', '0.05606'),
    @('E5', '0.36%'),
    @('D6', '6.468'),
    @('E6', '-0.04%'),
    @('D7', '0.8185'),
    @('E7', '-0.03%'),
    @('D8', '0.8283'),
    @('E8', '-1.40%'),
    @('B9', 'WazirX'),
    @('C9', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D9', '0.1332'),
    @('E9', '-0.40%'),
    @('B10', 'MandalaExchangeToken'),
    @('C10', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D10', '0.06941'),
    @('E10', '-0.21%'),
    @('B11', 'BitrueCoin'),
    @('C11', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D11', '0.02895'),
    @('E11', '1.45%'),
    @('B12', 'BitMartToken'),
    @('C12', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D12', '0.09385'),
    @('E12', '-0.05%'),
    @('B13', 'BitForexToken'),
    @('C13', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D13', '0.001522'),
    @('E13', '0.36%'),
    @('B14', 'TigerCash'),
    @('C14', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D14', '0.006247'),
    @('E14', '2.18%'),
    @('B15', 'LEO'),
    @('C15', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D15', '3.655'),
    @('E15', '3.71%'),
    @('B16', 'GateToken'),
    @('C16', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('D16', '3.034'),
    @('E16', '0.61%'),
    @('B17', 'BTSEToken'),
    @('C17', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'),
    @('D17', '2.182'),
    @('E17', '7.89%'),
    @('B18', 'BitpandaEcosystemToken'),
    @('C18', 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'),
    @('D18', '0.3112'),
    @('E18', '-2.10%'),
    @('B19', 'LiechtensteinCryptoassetsExchange'),
    @('C19', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D19', '0.03064'),
    @('E19', '-4.75%'),
    @('B20', 'ProBitToken'),
    @('C20', 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'),
    @('D20', '0.1299'),
    @('E20', '-2.23%'),
    @('B21', 'MCDex'),
    @('C21', 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'),
    @('D21', '3.758'),
    @('E21', '0.25%'),
    @('B22', 'CoinExToken'),
    @('C22', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'),
    @('D22', '0.04607'),
    @('E22', '-1.78%'),
    @('B23', 'ZBToken'),
    @('C23', 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'),
    @('D23', '0.1341'),
    @('E23', '-2.45%'),
    @('B24', 'One'),
    @('C24', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('D24', '0.0005947'),
    @('E24', '-93.93%'),
    @('D25', '0.001228'),
    @('E25', '-1.81%'),
    @('D26', '0.004487'),
    @('D27', '0.00009597'),
    @('E28', '0.63%'),
    @('D40', '0.03641'),
    @('E40', '-0.46%'),
    @('B41', 'BKEXToken'),
    @('C41', 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'),
    @('D41', '0.1370'),
    @('E41', '0.85%'),
    @('B42', 'CEJI'),
    @('C42', 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'),
    @('D42', '0.002569'),
    @('E42', '1.78%'),
    @('B43', 'KickToken'),
    @('C43', 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'),
    @('D43', '0.003394'),
    @('E43', '-45.09%'),
    @('D44', '0.009038'),
    @('E44', '7.37%'),
    @('D45', '0.00005345'),
    @('E45', '0.62%'),
    @('E46', '-0.04%'),
    @('E47', '8.22%'),
    @('D48', '0.002935'),
    @('E48', '38.28%'),
    @('E49', '-0.04%'),
    @('E50', '-0.04%')
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $rng = $ws.Range($cellRef)
    # Force text storage: without an explicit "@" (Text) number format, Excel
    # would parse numeric- or percent-looking strings (e.g. "244.14", "-0.40%")
    # into real numbers and silently reformat/round them.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    # Drop the "quote-prefixed text" style Excel stamps on text-that-looks-like-a-
    # number so the cell keeps the original (default) style, matching the source.
    $rng.Style = "Normal"
}
